$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" on every sheet ---
# "Overview" sheet: status shown in columns E (zh-cn) and F (de-de), rows 2-4.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# "zh-cn" / "de-de" sheets: status shown in column C, rows 2-4.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- 2. Column width changes: shrink the status-related columns ---
# width 17.2159881591797 -> 13.4101845877511 (nearest reachable ColumnWidth is 12.5)
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
